$wb = $excel.ActiveWorkbook

# --- ev_charging_uc: re-shuffled comma-separated timeslice lists ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")

$wsEv.Range("C13").Value = "WaD,RaP,FaD,SaD,RaD,FaP,SaP,WaP"
$wsEv.Range("C14").Value = "SaN,WaN,WaP,RaP,FaP,SaP,FaN,RaN"

# --- re_profiles: re-ordered season rows (M4:N7) ---
$wsRe = $wb.Worksheets.Item("re_profiles")

$wsRe.Range("M4").Value = "W"
$wsRe.Range("N4").Value = 0.29696276080640899

$wsRe.Range("M5").Value = "R"
$wsRe.Range("N5").Value = 0.34481908618716439

$wsRe.Range("M6").Value = "S"
$wsRe.Range("N6").Value = 0.27551721102209698

$wsRe.Range("M7").Value = "F"
$wsRe.Range("N7").Value = 0.2827009419843296
